# Adds a "Known Bugs:" section (preceded by a blank paragraph) at the end
# of the document, right after the "Fours is ..." bullet point.

$d = $word.ActiveDocument

# Locate the paragraph that ends the existing "To Do" list so we can anchor
# the new content right after it, regardless of exact paragraph index.
$anchor = $d.Content
$anchor.Find.Execute(
    "Fours is 20 minutes then 1 minute then 20 minutes then 1 minute",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0
) | Out-Null

# Collapse to the end of the found text (end of that paragraph's content).
$anchor.Collapse(0)

# --- New blank paragraph -------------------------------------------------
$anchor.InsertParagraphAfter()
$blankPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$blankPara.Range.ListFormat.RemoveNumbers()
$blankPara.Style = "Normal"

# --- "Known Bugs:" paragraph ---------------------------------------------
$r = $blankPara.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$bugsHeaderPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$bugsHeaderPara.Range.ListFormat.RemoveNumbers()
$bugsHeaderPara.Style = "Normal"
$bugsHeaderPara.Range.InsertAfter("Known Bugs:")

# --- Bug description paragraph -------------------------------------------
$r2 = $bugsHeaderPara.Range
$r2.Collapse(0)
$r2.InsertParagraphAfter()
$bugPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$bugPara.Range.ListFormat.RemoveNumbers()
$bugPara.Style = "Normal"
$bugPara.Range.InsertAfter("When refreshing page, chosen timer always defaults back to fours")
